$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("Y5").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-01-01-02-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_002.xlsx"
$ws.Range("V5").Value = 381.7622842221533
$ws.Range("X5").Value = 1238.6643104
$ws.Range("Z5").Value = 381.7622842221533
$ws.Range("AA5").Value = 363.7622842221533
$ws.Range("AH5").Value = 0.1600912605362297
$ws.Range("AI5").Value = 0.1600912605362297
$ws.Range("AJ5").Value = 0.1525429960566111
$ws.Range("AL5").Value = 100.423363368773
$ws.Range("AM5").Value = 263.5580617332523
$ws.Range("AQ5").Value = 560.7096415709138
$ws.Range("AR5").Value = 2384.654121302006
$ws.Range("AT5").Value = 1804.585617161007
$ws.Range("AU5").Value = 4.409614257598677
$ws.Range("AV5").Value = 264.0712239680002
$ws.Range("AW5").Value = 0.513162234747914
$ws.Range("AX5").Value = 167.3722722812883
$ws.Range("AY5").Value = 40.05415611619642
$ws.Range("AZ5").Value = 0.2191408798719717

# Row 9
$ws.Range("Y9").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-02-01-03-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_002.xlsx"
$ws.Range("V9").Value = 378.7312921699576
$ws.Range("X9").Value = 1201.0645585
$ws.Range("Z9").Value = 378.7312921699576
$ws.Range("AA9").Value = 360.7312921699576
$ws.Range("AH9").Value = 0.151280509043063
$ws.Range("AI9").Value = 0.151280509043063
$ws.Range("AJ9").Value = 0.144090585160161
$ws.Range("AL9").Value = 103.8987204790237
$ws.Range("AM9").Value = 257.1597640051947
$ws.Range("AQ9").Value = 606.4228340536064
$ws.Range("AR9").Value = 2503.503554857481
$ws.Range("AT9").Value = 1785.347397625189
$ws.Range("AU9").Value = 6.807237924860658
$ws.Range("AV9").Value = 258.0489082644798
$ws.Range("AW9").Value = 0.8891442592851279
$ws.Range("AX9").Value = 173.1645341317061
$ws.Range("AY9").Value = 54.01400200883849
$ws.Range("AZ9").Value = 0.3271923142607898

# Row 13
$ws.Range("Y13").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-03-01-04-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V13").Value = 319.8925843792476
$ws.Range("X13").Value = 1219.1102593
$ws.Range("Z13").Value = 319.8925843792476
$ws.Range("AA13").Value = 301.8925843792476
$ws.Range("AH13").Value = 0.1350252134828842
$ws.Range("AI13").Value = 0.1350252134828842
$ws.Range("AJ13").Value = 0.1274274948692808
$ws.Range("AL13").Value = 94.14891521222496
$ws.Range("AM13").Value = 208.2638731769835
$ws.Range("AT13").Value = 1469.804033549719
$ws.Range("AU13").Value = 22.61344177428888
$ws.Range("AV13").Value = 210.9446171702643
$ws.Range("AW13").Value = 2.680743993280815
$ws.Range("AX13").Value = 156.9148586870416
$ws.Range("AY13").Value = 95.08190639471097
$ws.Range("AZ13").Value = 0.5202040099608033

# Row 17
$ws.Range("Y17").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-04-01-05-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V17").Value = 328.3788769401734
$ws.Range("X17").Value = 1462.5820843
$ws.Range("Z17").Value = 328.3788769401734
$ws.Range("AA17").Value = 310.3788769401734
$ws.Range("AH17").Value = 0.1265334091252801
$ws.Range("AI17").Value = 0.1265334091252801
$ws.Range("AJ17").Value = 0.119597514266641
$ws.Range("AL17").Value = 93.06812538604369
$ws.Range("AM17").Value = 217.3107515541297
$ws.Range("AQ17").Value = 693.4528980543224
$ws.Range("AR17").Value = 2595.195049356863
$ws.Range("AT17").Value = 1527.011154796215
$ws.Range("AU17").Value = 6.115309273527662
$ws.Range("AV17").Value = 218.0670793198073
$ws.Range("AW17").Value = 0.7563277656775866
$ws.Range("AX17").Value = 155.1135423100728

# Row 21
$ws.Range("Y21").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-05-01-06-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V21").Value = 331.9143485514842
$ws.Range("X21").Value = 1511.2594947
$ws.Range("Z21").Value = 331.9143485514842
$ws.Range("AA21").Value = 313.9143485514842
$ws.Range("AH21").Value = 0.1295202985453993
$ws.Range("AI21").Value = 0.1295202985453993
$ws.Range("AJ21").Value = 0.1224963015895835
$ws.Range("AL21").Value = 89.44986790154844
$ws.Range("AM21").Value = 224.5182362881832
$ws.Range("AQ21").Value = 700.306823529412
$ws.Range("AR21").Value = 2562.643479663862
$ws.Range("AT21").Value = 1595.896016282618
$ws.Range("AU21").Value = 11.32613513684279
$ws.Range("AV21").Value = 225.8428845473684
$ws.Range("AW21").Value = 1.324648259185148
$ws.Range("AX21").Value = 149.0831131692474
$ws.Range("AY21").Value = 9.825354026807531
$ws.Range("AZ21").Value = 0.05375563824742705

# Row 25
$ws.Range("Y25").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-06-01-07-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V25").Value = 311.1362173684662
$ws.Range("X25").Value = 1202.5276552
$ws.Range("Z25").Value = 311.1362173684662
$ws.Range("AA25").Value = 293.1362173684662
$ws.Range("AH25").Value = 0.1399289434667034
$ws.Range("AI25").Value = 0.1399289434667034
$ws.Range("AJ25").Value = 0.1318337078695635
$ws.Range("AL25").Value = 96.84935292972744
$ws.Range("AM25").Value = 197.5573968366753
$ws.Range("AQ25").Value = 545.3665455761844
$ws.Range("AR25").Value = 2223.53009792075
$ws.Range("AT25").Value = 1427.043869798723
$ws.Range("AU25").Value = 17.24144434714563
$ws.Range("AV25").Value = 199.9453416704439
$ws.Range("AW25").Value = 2.38794483376855
$ws.Range("AX25").Value = 161.4155882162124
$ws.Range("AY25").Value = 224.7315877649988
$ws.Range("AZ25").Value = 1.270532397936547

# Row 29
$ws.Range("Y29").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-07-01-08-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V29").Value = 223.2638435420032
$ws.Range("X29").Value = 1169.5332397
$ws.Range("Z29").Value = 223.2638435420032
$ws.Range("AA29").Value = 205.2638435420032
$ws.Range("AH29").Value = 0.1039746882379047
$ws.Range("AI29").Value = 0.1039746882379047
$ws.Range("AJ29").Value = 0.09559203048826245
$ws.Range("AL29").Value = 64.44381977257362
$ws.Range("AM29").Value = 141.5501548826628
$ws.Range("AQ29").Value = 558.7355527353792
$ws.Range("AR29").Value = 2147.290338886641
$ws.Range("AT29").Value = 1033.1591836436
$ws.Range("AU29").Value = 31.73312956314376
$ws.Range("AV29").Value = 145.0147040711629
$ws.Range("AW29").Value = 3.46454918850011
$ws.Range("AX29").Value = 107.4063662876227
$ws.Range("AY29").Value = 133.4519858267323
$ws.Range("AZ29").Value = 0.7301311132331286

# Row 33
$ws.Range("Y33").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-08-01-09-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V33").Value = 233.3653695316492
$ws.Range("X33").Value = 1086.3745479
$ws.Range("Z33").Value = 233.3653695316492
$ws.Range("AA33").Value = 215.3653695316492
$ws.Range("AH33").Value = 0.1087443545646453
$ws.Range("AI33").Value = 0.1087443545646453
$ws.Range("AJ33").Value = 0.1003566559695538
$ws.Range("AL33").Value = 65.956972776347
$ws.Range("AM33").Value = 149.9350775801196
$ws.Range("AQ33").Value = 530.2951205625106
$ws.Range("AT33").Value = 1071.905169927031
$ws.Range("AU33").Value = 15.61085430927645
$ws.Range("AV33").Value = 151.7337745682448
$ws.Range("AW33").Value = 1.798696988125188
$ws.Range("AX33").Value = 109.9282879605784
$ws.Range("AY33").Value = 96.26572638098139
$ws.Range("AZ33").Value = 0.5266808248173875

# Row 37
$ws.Range("Y37").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-09-01-10-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V37").Value = 284.9811300356272
$ws.Range("X37").Value = 1218.0763376
$ws.Range("Z37").Value = 284.9811300356272
$ws.Range("AA37").Value = 266.9811300356272
$ws.Range("AH37").Value = 0.1274697215339408
$ws.Range("AI37").Value = 0.1274697215339408
$ws.Range("AJ37").Value = 0.1194184692025177
$ws.Range("AL37").Value = 87.03783460116286
$ws.Range("AM37").Value = 180.1116038895606
$ws.Range("AQ37").Value = 594.7664934678277
$ws.Range("AT37").Value = 1275.759182717254
$ws.Range("AU37").Value = 7.246192519998683
$ws.Range("AV37").Value = 181.0260442021478
$ws.Range("AW37").Value = 0.9144403125872189
$ws.Range("AX37").Value = 145.0630576686048
$ws.Range("AY37").Value = 29.77037532414831
$ws.Range("AZ37").Value = 0.1683084550962807

# Row 41
$ws.Range("Y41").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-10-01-11-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V41").Value = 359.5342253400781
$ws.Range("X41").Value = 1603.303729
$ws.Range("Z41").Value = 359.5342253400781
$ws.Range("AA41").Value = 341.5342253400781
$ws.Range("AH41").Value = 0.1334921997668373
$ws.Range("AI41").Value = 0.1334921997668373
$ws.Range("AJ41").Value = 0.1268089428570668
$ws.Range("AL41").Value = 83.26105551668992
$ws.Range("AM41").Value = 258.5388161535467
$ws.Range("AQ41").Value = 832.0349031040995
$ws.Range("AT41").Value = 1784.609881155756
$ws.Range("AU41").Value = 1.908148349886507
$ws.Range("AV41").Value = 258.8310040101808
$ws.Range("AW41").Value = 0.2921878566340159
$ws.Range("AX41").Value = 138.7684258611499
$ws.Range("AY41").Value = 48.55433448145752
$ws.Range("AZ41").Value = 0.2656463301585349

# Row 45
$ws.Range("Y45").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-11-01-12-01-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V45").Value = 396.2144451703995
$ws.Range("X45").Value = 1358.5385934
$ws.Range("Z45").Value = 396.2144451703995
$ws.Range("AA45").Value = 378.2144451703995
$ws.Range("AH45").Value = 0.1575011201017935
$ws.Range("AI45").Value = 0.1575011201017935
$ws.Range("AJ45").Value = 0.1503458530579253
$ws.Range("AL45").Value = 102.2603489709766
$ws.Range("AM45").Value = 276.0866970513032
$ws.Range("AQ45").Value = 754.4330897207806
$ws.Range("AT45").Value = 1920.633769695537
$ws.Range("AU45").Value = 9.55135447016875
$ws.Range("AV45").Value = 277.4947509882299
$ws.Range("AW45").Value = 1.408053936926675
$ws.Range("AX45").Value = 170.4339149516277
$ws.Range("AY45").Value = 23.45441960430264
$ws.Range("AZ45").Value = 0.1326008518803594

# Row 49
$ws.Range("Y49").Value = "MPC-MPC-optimal-Simple-unconscious-1.5-0.6-0.6-flex-438-unconscious-12-01-12-31-Sum-ALL-Sum-nan-nan-minimize_cap-2023-09-08_001.xlsx"
$ws.Range("V49").Value = 353.5612346795531
$ws.Range("X49").Value = 1102.9779828
$ws.Range("Z49").Value = 353.5612346795531
$ws.Range("AA49").Value = 335.5612346795531
$ws.Range("AH49").Value = 0.1713214432365774
$ws.Range("AI49").Value = 0.1713214432365774
$ws.Range("AJ49").Value = 0.1625993728403323
$ws.Range("AM49").Value = 229.4392872459554
$ws.Range("AT49").Value = 1634.710088034392
$ws.Range("AU49").Value = 22.84944035846
$ws.Range("AV49").Value = 232.4805308339622
$ws.Range("AW49").Value = 3.041243588006854
$ws.Range("AY49").Value = 85.15572851986794
$ws.Range("AZ49").Value = 0.4814325971279066
